$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    0.5656572556951005,0.1011745124650645,0.05147103804321063,0,1.202792717630317,1.074386414929208,1.065747477021972,1.091568419121032,0,0.3442417395016832,0.3037225118552982,
    0.5252628992619748,0.09935451473029389,0.05065505260682457,0,1.196938372930227,1.070033925374418,1.068421079425974,1.093844930269817,0,0.3073039585933657,0.2926883187876967,
    0.5006870167006809,0.09821731973729442,0.05014546337182324,0,1.194012476497484,1.067991152768897,1.070520786456527,1.095767598040105,0,0.2846893843581881,0.2860901357356767,
    0.4907295082801966,0.09774894952183644,0.04993565995045657,0,1.192988189999326,1.067316834754081,1.071491638621268,1.096683015764668,0,0.2754905014299567,0.2834457935115466,
    0.4890795489263269,0.09767087785827044,0.04990069330305147,0,1.192828254050958,1.067214409246915,1.071659806277168,1.096842985469173,0,0.2739640539190589,0.2830093884369944,
    0.5005524934260563,0.09821102318822028,0.05014264254300116,0,1.193997982371421,1.067981418703837,1.070533413239687,1.09577940966183,0,0.2845652568251751,0.2860542931502579,
    0.5516825319592726,0.1005510685447035,0.05119147300122862,0,1.200635281624002,1.07275485848352,1.066574260115999,1.092244376802363,0,0.3314921927980947,0.2998812107112201,
    0.6537316979938339,0.1049836078914979,0.05317971398915233,0,1.218963978473354,1.087123031628678,1.062445586976835,1.089481021723572,0,0.4240260843091335,0.3284008993743157,
    0.729785019388288,0.1081453303912383,0.05459815773000543,0,1.235682227272761,1.100750404634127,1.061630082472277,1.089999337748132,0,0.4923177158918577,0.3502161561199273,
    0.764616265387275,0.1095631867557287,0.05523415038083357,0,1.243997015044044,1.10762082765676,1.061741132593767,1.090790157955794,0,0.5234516165358798,0.3603289084425825,
    0.77783933075969,0.1100971571192204,0.05547364073434125,0,1.247247825025681,1.110319282948524,1.061852518868662,1.09116954344087,0,0.5352507580993802,0.3641855526687721,
    0.774990036333719,0.1099822879103414,0.05542212229864418,0,1.246543158138479,1.109733814267585,1.061825445883841,1.091084279886402,0,0.5327091885300206,0.3633537471339849,
    0.765703471554076,0.1096071757829478,0.05525388045648327,0,1.244262412466327,1.107840890258004,1.061748906742693,1.090819767927009,0,0.5244221522952444,0.3606456524499038,
    0.7600194964959712,0.1093770256531457,0.05515065173540989,0,1.24287870105951,1.106694029993974,1.061711054215223,1.090668157664332,0,0.5193473190492739,0.3589904039667289,
    0.7275133898390322,0.1080522587800914,0.05455640648478521,0,1.23515312713117,1.100314935107548,1.061632526326349,1.089958833550568,0,0.4902843765839009,0.3495590650563543,
    0.7076316111472067,0.1072343226867005,0.05418947334470658,0,1.230595579454132,1.096573670505833,1.061707824189313,1.089665919273912,0,0.472472308652641,0.3438216238274094,
    0.6962182094396496,0.1067619449989721,0.05397755190961817,0,1.228040978286771,1.094484961666666,1.061796503031502,1.089549678991531,0,0.4622336693346654,0.3405393777143075,
    0.6923576323197267,0.1066016759576556,0.05390564975547107,0,1.227187499291773,1.093788602262549,1.061834319850163,1.089519290444301,0,0.4587681496107621,0.3394311189208707,
    0.7097457793824447,0.1073215923763726,0.05422862427151642,0,1.231073825681776,1.09696539509406,1.061695112889922,1.089691693473199,0,0.4743677746092203,0.3444305442279756,
    0.7684302625029602,0.1097174351561776,0.0553033338105493,0,1.244929548462721,1.108394259497246,1.061769506284591,1.090895291680297,0,0.5268560029378193,0.3614403482435478,
    0.8069774388156361,0.1112661124971694,0.05599786524703632,0,1.25458071594187,1.116427877790827,1.062222267467916,1.092147790967594,0,0.5612147886716627,0.3727155749265876,
    0.7863865304301214,0.1104411241488421,0.05562790390894889,0,1.249375156448281,1.112088478760896,1.061943634983137,1.091436645759025,0,0.5428719603580987,0.3666832853453315,
    0.7087899112572131,0.1072821444013172,0.05421092714679787,0,1.230857406248475,1.09678810269196,1.061700718279283,1.089679878485164,0,0.4735108289594905,0.344155200591203,
    0.6259347594419751,0.1038011778634242,0.05264923460156012,0,1.213435420969873,1.082698072909295,1.063173153093686,1.089781544329739,0,0.3989390192569431,0.3205345975156035
)

$idx = 0
for ($r = 2; $r -le 25; $r++) {
  for ($c = 2; $c -le 12; $c++) {
    $ws.Cells.Item($r, $c).Value = $data[$idx]
    $idx = $idx + 1
  }
}
